# Weekly update: a new "Arándano (blue)" price record for Vega Modelo de
# Temuco is inserted at the top of the data (row 19), pushing every
# existing weekly record down by one row (old row 67 becomes row 68).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("19:19").Insert()

$ws.Cells.Item(19, 1).Value2  = 10
$ws.Cells.Item(19, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value2  = "La Araucanía"
$ws.Cells.Item(19, 4).Value2  = 44536
$ws.Cells.Item(19, 5).Value2  = 9
$ws.Cells.Item(19, 6).Value2  = "Fruta"
$ws.Cells.Item(19, 7).Value2  = 100101
$ws.Cells.Item(19, 8).Value2  = "Berries"
$ws.Cells.Item(19, 9).Value2  = 100101001
$ws.Cells.Item(19, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(19, 11).Value2 = "Sin especificar"
$ws.Cells.Item(19, 12).Value2 = "Primera"
$ws.Cells.Item(19, 13).Value2 = 600
$ws.Cells.Item(19, 14).Value2 = 3000
$ws.Cells.Item(19, 15).Value2 = 3000
$ws.Cells.Item(19, 16).Value2 = 3000
$ws.Cells.Item(19, 17).Value2 = "`$/kilo"
$ws.Cells.Item(19, 18).Value2 = "Región del Maule"
$ws.Cells.Item(19, 19).Value2 = 3000
$ws.Cells.Item(19, 20).Value2 = 1
